# Apply the cryptos-list refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.116.10'
$ws.Cells.Item(2, 5).Value = '  -1.00%  '

$ws.Cells.Item(3, 4).Value = '1.794.65'
$ws.Cells.Item(3, 5).Value = '  -0.52%  '

$ws.Cells.Item(4, 4).Value = '0.9994'
$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).Value = '''316.80'
$ws.Cells.Item(5, 5).Value = '  +0.27%  '

$ws.Cells.Item(6, 4).Value = '0.9992'
$ws.Cells.Item(6, 5).Value = '  -0.08%  '

$ws.Cells.Item(7, 4).Value = '0.5343'
$ws.Cells.Item(7, 5).Value = '  -3.18%  '

$ws.Cells.Item(8, 4).Value = '0.3765'
$ws.Cells.Item(8, 5).Value = '  -2.36%  '

$ws.Cells.Item(9, 4).Value = '0.07469'
$ws.Cells.Item(9, 5).Value = '  -1.62%  '

$ws.Cells.Item(10, 4).Value = '41.83'
$ws.Cells.Item(10, 5).Value = '  -1.69%  '

$ws.Cells.Item(11, 4).Value = '1.096'
$ws.Cells.Item(11, 5).Value = '  -2.69%  '

$ws.Cells.Item(12, 4).Value = '0.9997'
$ws.Cells.Item(12, 5).Value = '  -0.07%  '

$ws.Cells.Item(13, 4).Value = '20.65'
$ws.Cells.Item(13, 5).Value = '  -2.61%  '

$ws.Cells.Item(14, 4).Value = '6.112'
$ws.Cells.Item(14, 5).Value = '  -1.26%  '

$ws.Cells.Item(15, 4).Value = '7.234'
$ws.Cells.Item(15, 5).Value = '  -2.35%  '

$ws.Cells.Item(16, 4).Value = '1.782.96'
$ws.Cells.Item(16, 5).Value = '  -1.47%  '

$ws.Cells.Item(17, 4).Value = '89.16'
$ws.Cells.Item(17, 5).Value = '  -3.20%  '

$ws.Cells.Item(18, 4).Value = '0.00001057'
$ws.Cells.Item(18, 5).Value = '  -1.28%  '

$ws.Cells.Item(19, 4).Value = '0.06467'
$ws.Cells.Item(19, 5).Value = '  +0.40%  '

$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20, 4).Value = '0.9982'
$ws.Cells.Item(20, 5).Value = '  -0.15%  '

$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = '17.35'
$ws.Cells.Item(21, 5).Value = '  +0.14%  '

$ws.Cells.Item(22, 4).Value = '5.912'
$ws.Cells.Item(22, 5).Value = '  -1.10%  '

$ws.Cells.Item(23, 4).Value = '28.141.09'
$ws.Cells.Item(23, 5).Value = '  -1.01%  '

$ws.Cells.Item(24, 4).Value = '11.21'
$ws.Cells.Item(24, 5).Value = '  -2.08%  '

$ws.Cells.Item(25, 4).Value = '2.093'
$ws.Cells.Item(25, 5).Value = '  -2.08%  '

$ws.Cells.Item(26, 4).Value = '154.81'
$ws.Cells.Item(26, 5).Value = '  -2.77%  '

$ws.Cells.Item(27, 4).Value = '20.25'
$ws.Cells.Item(27, 5).Value = '  -2.08%  '

$ws.Cells.Item(28, 4).Value = '1.992.05'
$ws.Cells.Item(28, 5).Value = '  -1.18%  '

$ws.Cells.Item(29, 4).Value = '2.292'
$ws.Cells.Item(29, 5).Value = '  -4.81%  '

$ws.Cells.Item(30, 4).Value = '120.45'
$ws.Cells.Item(30, 5).Value = '  -2.80%  '

$ws.Cells.Item(31, 4).Value = '1.115'
$ws.Cells.Item(31, 5).Value = '  -0.43%  '

$ws.Cells.Item(32, 4).Value = '0.1046'
$ws.Cells.Item(32, 5).Value = '  +2.60%  '

$ws.Cells.Item(33, 4).Value = '3.654'
$ws.Cells.Item(33, 5).Value = '  -0.74%  '

$ws.Cells.Item(34, 5).Value = '  -3.14%  '

$ws.Cells.Item(35, 4).Value = '0.2262'
$ws.Cells.Item(35, 5).Value = '  -2.28%  '

$ws.Cells.Item(36, 4).Value = '0.06517'
$ws.Cells.Item(36, 5).Value = '  +1.00%  '

$ws.Cells.Item(37, 4).Value = '0.02287'
$ws.Cells.Item(37, 5).Value = '  -1.62%  '

$ws.Cells.Item(38, 4).Value = '5.023'
$ws.Cells.Item(38, 5).Value = '  -2.23%  '

$ws.Cells.Item(39, 4).Value = '8.485'
$ws.Cells.Item(39, 5).Value = '  -3.80%  '

$ws.Cells.Item(40, 5).Value = '  +4.49%  '

$ws.Cells.Item(41, 4).Value = '0.6163'
$ws.Cells.Item(41, 5).Value = '  -3.94%  '

$ws.Cells.Item(42, 5).Value = '  -4.77%  '

$ws.Cells.Item(43, 4).Value = '1.173'
$ws.Cells.Item(43, 5).Value = '  +1.06%  '

$ws.Cells.Item(44, 4).Value = '0.9985'
$ws.Cells.Item(44, 5).Value = '  -0.11%  '

$ws.Cells.Item(45, 5).Value = '  -2.43%  '

$ws.Cells.Item(46, 4).Value = '3.672'
$ws.Cells.Item(46, 5).Value = '  -0.16%  '

$ws.Cells.Item(47, 4).Value = '0.5776'
$ws.Cells.Item(47, 5).Value = '  -3.54%  '

$ws.Cells.Item(48, 4).Value = '127.27'
$ws.Cells.Item(48, 5).Value = '  +0.02%  '

$ws.Cells.Item(49, 4).Value = '1.189'
$ws.Cells.Item(49, 5).Value = '  +3.50%  '

$ws.Cells.Item(50, 4).Value = '''1.930'
$ws.Cells.Item(50, 5).Value = '  -2.78%  '

$ws.Cells.Item(51, 4).Value = '0.06815'
$ws.Cells.Item(51, 5).Value = '  -1.31%  '
